# feat: add 2022-Q1 data
#
# 1) Insert a new worksheet "2022-Q1" (fund-holding detail) right after "2021-Q4",
#    copying the column layout/styling from "2021-Q4".
# 2) Update the "总计" (totals) sheet with a new leading row for 2022-Q1,
#    pushing the existing 2021-Q4 / 2021-Q3 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q1" sheet
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $template)
$q1.Name = "2022-Q1"

# Copy header-row formatting (bold/border/alignment style used across all sheets).
$template.Range("A1:H1").Copy()
$q1.Range("A1:H1").PasteSpecial(-4122)

# Copy the data-row formatting pattern (index column style) and let Excel tile it
# down across all 10 data rows (rows 2-11).
$template.Range("A2:H2").Copy()
$q1.Range("A2:H11").PasteSpecial(-4122)

# Columns B-G hold text (fund code keeps leading zeros, numeric-looking figures
# are stored as text to match the other quarters' sheets) - force text format
# before writing so Excel doesn't reinterpret them as numbers.
$q1.Range("B2:G11").NumberFormat = "@"

# Header row
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

$q1Rows = @(
    @("010846","南方卓越优选3个月持有期混合A","26.01","60.36","2.01","0.5228",10),
    @("005123","南方优享分红灵活配置混合A","4.86","93.90","7.14","0.3470",5),
    @("010847","南方卓越优选3个月持有期混合C","3.42","60.36","2.01","0.0687",10),
    @("004703","南方兴盛先锋灵活配置混合","1.09","53.08","6.19","0.0675",2),
    @("000066","诺安鸿鑫混合","0.74","81.34","4.92","0.0364",4),
    @("001731","广发百发大数据策略价值灵活配置混合A","0.24","88.87","2.88","0.0069",7),
    @("001732","广发百发大数据策略价值灵活配置混合E","0.24","88.87","2.88","0.0069",7),
    @("006587","南方优享分红灵活配置混合C","0.09","93.90","7.14","0.0064",5),
    @("005536","渤海汇金量化成长混合","0.61","88.57","0.69","0.0042",9),
    @("002952","建信多因子量化股票","0.10","91.47","2.95","0.0030",5)
)

for ($i = 0; $i -lt $q1Rows.Length; $i++) {
    $r = $i + 2
    $row = $q1Rows[$i]
    $q1.Cells.Item($r,1).Value = $i
    $q1.Cells.Item($r,2).Value = $row[0]
    $q1.Cells.Item($r,3).Value = $row[1]
    $q1.Cells.Item($r,4).Value = $row[2]
    $q1.Cells.Item($r,5).Value = $row[3]
    $q1.Cells.Item($r,6).Value = $row[4]
    $q1.Cells.Item($r,7).Value = $row[5]
    $q1.Cells.Item($r,8).Value = $row[6]
}

# ---------------------------------------------------------------------------
# 2) Update "总计" sheet: add a 2022-Q1 row on top, shift the rest down
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Push the existing two data rows (2021-Q4, 2021-Q3) down one row, preserving
# their formatting, then write the new 2022-Q1 summary into row 2.
$total.Range("A2:D2").Copy()
$total.Range("A3:A4").PasteSpecial(-4122)

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2021-Q3"
$total.Cells.Item(4,3).Value = 6
$total.Cells.Item(4,4).Value = 1.49

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q4"
$total.Cells.Item(3,3).Value = 3
$total.Cells.Item(3,4).Value = 0.54

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 10
$total.Cells.Item(2,4).Value = 1.07

# Restore the originally-active tab (adding a sheet shifts focus to it by default).
$wb.Worksheets.Item("2021-Q3").Activate()

Write-Output "2022-Q1 sheet added and 总计 sheet updated"
